{"js": "// Apply the \"Added many more features\" edits to the John Hunter Tomb of\n// the Scarab Queen review document.\n//\n// Each entry is an exact (old -> new) full-paragraph text swap. We use\n// Body.search with matchCase so we target the precise text runs described\n// in the diff without disturbing any other (similar) text elsewhere in\n// the document.\nconst replacements = [\n  {\n    from: \"Play John Hunter Tomb of the Scarab Queen for Free\",\n    to: \"Play John Hunter Tomb of the Scarab Queen Free\",\n  },\n  {\n    from: \"Enjoyable gameplay with multiple free spins and multipliers\",\n    to: \"Original and linear gameplay\",\n  },\n  {\n    from: \"Themed symbols and design with Ancient Egypt artifacts\",\n    to: \"Free spins and multipliers up to x25\",\n  },\n  {\n    from: \"Medium volatility and high maximum win of x10,500\",\n    to: \"Themed symbols and immersive design\",\n  },\n  {\n    from: \"Available to play on both desktop and mobile devices\",\n    to: \"Available on both desktop and mobile devices\",\n  },\n  {\n    from: \"Limited paylines compared to other slot games\",\n    to: \"Medium volatility may not appeal to all players\",\n  },\n  {\n    from: \"Cards from J to Ace lack creativity in design\",\n    to: \"Maximum bet amount of \\u20ac125 per spin may be too high for some players\",\n  },\n  {\n    from:\n      \"Read our review of John Hunter Tomb of the Scarab Queen, a Pragmatic Play slot game. Play for free with enjoyable gameplay and multipliers up to x25.\",\n    to: \"Read our review of John Hunter Tomb of the Scarab Queen and play for free.\",\n  },\n];\n\nfor (const { from, to } of replacements) {\n  const results = context.document.body.search(from, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edits to the John Hunter Tomb of\n# the Scarab Queen review document via Word COM interop Find/Replace.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ From = \"Play John Hunter Tomb of the Scarab Queen for Free\"; To = \"Play John Hunter Tomb of the Scarab Queen Free\" },\n    @{ From = \"Enjoyable gameplay with multiple free spins and multipliers\"; To = \"Original and linear gameplay\" },\n    @{ From = \"Themed symbols and design with Ancient Egypt artifacts\"; To = \"Free spins and multipliers up to x25\" },\n    @{ From = \"Medium volatility and high maximum win of x10,500\"; To = \"Themed symbols and immersive design\" },\n    @{ From = \"Available to play on both desktop and mobile devices\"; To = \"Available on both desktop and mobile devices\" },\n    @{ From = \"Limited paylines compared to other slot games\"; To = \"Medium volatility may not appeal to all players\" },\n    @{ From = \"Cards from J to Ace lack creativity in design\"; To = \"Maximum bet amount of \u20ac125 per spin may be too high for some players\" },\n    @{ From = \"Read our review of John Hunter Tomb of the Scarab Queen, a Pragmatic Play slot game. Play for free with enjoyable gameplay and multipliers up to x25.\"; To = \"Read our review of John Hunter Tomb of the Scarab Queen and play for free.\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.From\n    $find.Replacement.Text = $r.To\n    $find.Execute(\n        [ref]$find.Text,\n        [ref]$true,\n        [ref]$false,\n        [ref]$false,\n        [ref]$false,\n        [ref]$false,\n        [ref]$true,\n        [ref]1,\n        [ref]$false,\n        [ref]$find.Replacement.Text,\n        [ref]2\n    )\n}\n"}
